# The edit removes the last column (AJ) from the worksheet, which held a
# per-row "conv_attn_<timestamp>" run identifier. After the deletion the
# data that used to sit in column AI ("nan" for every data row, with the
# "loss_tr" header in row 1) becomes the last populated column, and the
# sheet's used range shrinks from A1:AJ12 to A1:AI12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the entire AJ column (shifts nothing else, it's the last one).
$ws.Range("AJ1:AJ12").EntireColumn.Delete() | Out-Null

# Match the cursor/selection position left behind in the saved workbook.
$ws.Range("AI8").Select() | Out-Null
